$d = $word.ActiveDocument
$d.Content.Find.Execute("Next.js, React, Django, Tailwind CSS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Next.js, Django, Tailwind CSS", 2)
